$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.065.05"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.497.16"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "320.44"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "107.51"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("D7").Value = "0.524"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "39.55"
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("D11").Value = "20.18"
$ws.Range("E11").Value = "  +7.86%  "
$ws.Range("D12").Value = "0.0813"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "2.887.84"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "2.497.27"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "0.836"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "47.924.80"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "0.0₃0939"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").Value = "277.56"
$ws.Range("E23").Value = "  +11.77%  "
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "25.60"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "34.95"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -4.87%  "
$ws.Range("D32").Value = "49.41"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "19.52"
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "0.0779"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("D38").Value = "4.62"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").Value = "120.63"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "21.16"
$ws.Range("E43").Value = "  -5.52%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "2.007.89"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "3.15"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.00"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.84"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "8.97"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "80.15"
$ws.Range("E51").Value = "  +2.54%  "
